# Updates the cryptos price table to the latest scraped values.
# (Matches the commit "Updated cryptos list ... with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    # Some of the replacement values look like plain numbers (e.g. "0.9997",
    # "306.41", "0.07560"). Assigning such a string straight to .Value lets the
    # COM layer coerce it into a real number, which silently drops significant
    # trailing zeros / reformats it (e.g. "0.07560" -> 0.0756). Forcing the cell
    # to Text format before the write keeps the exact original string, and
    # restoring the default style afterwards avoids leaving stray formatting on
    # the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '26.933.90'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '1.875.09'
$ws.Range("E3").Value = '  +0.94%  '
Set-TextValue "D4" '0.9997'
$ws.Range("E4").Value = '  -0.06%  '
Set-TextValue "D5" '306.41'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("E6").Value = '  -0.14%  '
Set-TextValue "D7" '0.5164'
$ws.Range("E7").Value = '  +1.68%  '
Set-TextValue "D8" '0.3714'
$ws.Range("E8").Value = '  +1.66%  '
Set-TextValue "D9" '0.07190'
$ws.Range("E9").Value = '  +0.84%  '
Set-TextValue "D10" '0.8981'
$ws.Range("E10").Value = '  +1.19%  '
Set-TextValue "D11" '20.68'
$ws.Range("E11").Value = '  -0.28%  '
Set-TextValue "D12" '0.07560'
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.896.20'
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D14" '94.89'
$ws.Range("E14").Value = '  +4.86%  '
Set-TextValue "D15" '5.251'
$ws.Range("E15").Value = '  +0.24%  '
Set-TextValue "D16" '0.9999'
$ws.Range("E16").Value = '  -0.08%  '
Set-TextValue "D17" '0.000008486'
$ws.Range("E17").Value = '  -0.52%  '
Set-TextValue "D18" '14.23'
$ws.Range("E18").Value = '  +1.40%  '
Set-TextValue "D19" '0.9992'
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = '26.951.91'
$ws.Range("E20").Value = '  +0.34%  '
Set-TextValue "D21" '5.028'
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").Value = '2.102.10'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").Value = '  +1.14%  '
Set-TextValue "D24" '6.426'
$ws.Range("E24").Value = '  -0.25%  '
Set-TextValue "D25" '145.96'
$ws.Range("E25").Value = '  +0.12%  '
Set-TextValue "D26" '1.784'
$ws.Range("E26").Value = '  -1.90%  '
$ws.Range("E27").Value = '  +1.15%  '
$ws.Range("E28").Value = '  +3.21%  '
Set-TextValue "D29" '114.52'
$ws.Range("E29").Value = '  +1.66%  '
Set-TextValue "D30" '4.903'
$ws.Range("E30").Value = '  +5.11%  '
Set-TextValue "D31" '4.739'
$ws.Range("E31").Value = '  +2.62%  '
Set-TextValue "D32" '0.09176'
$ws.Range("E32").Value = '  -0.37%  '
Set-TextValue "D33" '0.05029'
$ws.Range("E33").Value = '  -1.50%  '
Set-TextValue "D34" '0.7527'
$ws.Range("E34").Value = '  +2.98%  '
Set-TextValue "D35" '2.987'
$ws.Range("E35").Value = '  -2.67%  '
Set-TextValue "D36" '1.170'
$ws.Range("E36").Value = '  +1.87%  '
Set-TextValue "D37" '3.279'
$ws.Range("E37").Value = '  +2.79%  '
Set-TextValue "D38" '0.01990'
$ws.Range("E38").Value = '  -1.20%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D39" '2.487'
$ws.Range("E39").Value = '  +0.88%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D40" '0.5564'
$ws.Range("E40").Value = '  +4.77%  '
$ws.Range("E41").Value = '  -0.25%  '
Set-TextValue "D42" '6.572'
$ws.Range("E42").Value = '  +1.72%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D43" '115.95'
$ws.Range("E43").Value = '  -1.73%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D44" '8.718'
$ws.Range("E44").Value = '  +3.88%  '
Set-TextValue "D45" '0.1499'
$ws.Range("E45").Value = '  +2.05%  '
Set-TextValue "D46" '0.4761'
$ws.Range("E46").Value = '  +2.90%  '
$ws.Range("E47").Value = '  -0.18%  '
Set-TextValue "D48" '10.06'
$ws.Range("E48").Value = '  +1.53%  '
Set-TextValue "D49" '1.561'
$ws.Range("E49").Value = '  +0.43%  '
Set-TextValue "D50" '37.07'
$ws.Range("E50").Value = '  +0.44%  '
Set-TextValue "D51" '63.35'
$ws.Range("E51").Value = '  +0.34%  '
